# Update the build/version strings throughout the workbook.
$wb = $excel.ActiveWorkbook

$oldBuild = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newBuild = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A2").Value = "Version: " + $newBuild
$aboutSheet.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for KWK Pni' + [char]0x00F3 + 'wek Coal Mine, Poland, M1289, version ''' + $newBuild + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
$lastRow = $dataSheet.Cells.Item($dataSheet.Rows.Count, 19).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)
    if ($cell.Value2 -eq $oldBuild) {
        $cell.Value = $newBuild
    }
}
